$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns remain plain text, matching the source data
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.804.56'
$ws.Cells.Item(2, 5).Value = '  +0.75%  '
$ws.Cells.Item(3, 4).Value = '1.732.03'
$ws.Cells.Item(4, 4).Value = '0.9960'
$ws.Cells.Item(4, 5).Value = '  -0.41%  '
$ws.Cells.Item(5, 4).Value = '242.27'
$ws.Cells.Item(5, 5).Value = '  -1.51%  '
$ws.Cells.Item(6, 4).Value = '0.9968'
$ws.Cells.Item(6, 5).Value = '  -0.38%  '
$ws.Cells.Item(7, 4).Value = '0.4961'
$ws.Cells.Item(7, 5).Value = '  +0.86%  '
$ws.Cells.Item(8, 5).Value = '  -2.06%  '
$ws.Cells.Item(9, 4).Value = '0.06228'
$ws.Cells.Item(9, 5).Value = '  -0.71%  '
$ws.Cells.Item(10, 4).Value = '1.730.65'
$ws.Cells.Item(10, 5).Value = '  -0.53%  '
$ws.Cells.Item(11, 5).Value = '  +0.19%  '
$ws.Cells.Item(12, 4).Value = '0.06988'
$ws.Cells.Item(12, 5).Value = '  -0.80%  '
$ws.Cells.Item(13, 4).Value = '0.6143'
$ws.Cells.Item(13, 5).Value = '  +0.11%  '
$ws.Cells.Item(14, 4).Value = '4.506'
$ws.Cells.Item(14, 5).Value = '  -1.55%  '
$ws.Cells.Item(15, 4).Value = '77.21'
$ws.Cells.Item(15, 5).Value = '  -0.95%  '
$ws.Cells.Item(16, 4).Value = '0.9966'
$ws.Cells.Item(16, 5).Value = '  -0.39%  '
$ws.Cells.Item(17, 4).Value = '26.559.21'
$ws.Cells.Item(17, 5).Value = '  -0.22%  '
$ws.Cells.Item(18, 4).Value = '0.9961'
$ws.Cells.Item(18, 5).Value = '  -0.44%  '
$ws.Cells.Item(19, 4).Value = '0.000007183'
$ws.Cells.Item(19, 5).Value = '  -0.97%  '
$ws.Cells.Item(20, 4).Value = '11.42'
$ws.Cells.Item(20, 5).Value = '  -1.18%  '
$ws.Cells.Item(21, 4).Value = '1.951.32'
$ws.Cells.Item(21, 5).Value = '  -1.15%  '
$ws.Cells.Item(22, 4).Value = '4.437'
$ws.Cells.Item(22, 5).Value = '  -2.77%  '
$ws.Cells.Item(23, 4).Value = '8.538'
$ws.Cells.Item(23, 5).Value = '  -1.98%  '
$ws.Cells.Item(24, 4).Value = '5.124'
$ws.Cells.Item(24, 5).Value = '  -2.86%  '
$ws.Cells.Item(25, 4).Value = '138.50'
$ws.Cells.Item(25, 5).Value = '  -0.43%  '
$ws.Cells.Item(26, 4).Value = '15.36'
$ws.Cells.Item(26, 5).Value = '  -0.40%  '
$ws.Cells.Item(27, 4).Value = '1.415'
$ws.Cells.Item(27, 5).Value = '  -0.44%  '
$ws.Cells.Item(28, 5).Value = '  +0.04%  '
$ws.Cells.Item(29, 4).Value = '106.56'
$ws.Cells.Item(29, 5).Value = '  -0.81%  '
$ws.Cells.Item(30, 4).Value = '3.949'
$ws.Cells.Item(30, 5).Value = '  -1.76%  '
$ws.Cells.Item(31, 4).Value = '0.07988'
$ws.Cells.Item(31, 5).Value = '  -0.82%  '
$ws.Cells.Item(32, 4).Value = '3.661'
$ws.Cells.Item(32, 5).Value = '  -1.69%  '
$ws.Cells.Item(33, 4).Value = '0.04534'
$ws.Cells.Item(33, 5).Value = '  -1.73%  '
$ws.Cells.Item(34, 5).Value = '  -0.20%  '
$ws.Cells.Item(35, 4).Value = '1.004'
$ws.Cells.Item(35, 5).Value = '  -0.86%  '
$ws.Cells.Item(36, 4).Value = '0.6287'
$ws.Cells.Item(36, 5).Value = '  -1.54%  '
$ws.Cells.Item(37, 4).Value = '0.9445'
$ws.Cells.Item(37, 5).Value = '  +4.33%  '
$ws.Cells.Item(38, 4).Value = '2.021'
$ws.Cells.Item(38, 5).Value = '  -2.03%  '
$ws.Cells.Item(39, 4).Value = '2.423'
$ws.Cells.Item(39, 5).Value = '  -0.14%  '
$ws.Cells.Item(40, 4).Value = '0.9969'
$ws.Cells.Item(40, 5).Value = '  -0.61%  '
$ws.Cells.Item(41, 5).Value = '  -0.07%  '
$ws.Cells.Item(42, 4).Value = '99.91'
$ws.Cells.Item(42, 5).Value = '  -2.13%  '
$ws.Cells.Item(43, 4).Value = '5.498'
$ws.Cells.Item(43, 5).Value = '  +1.25%  '
$ws.Cells.Item(44, 4).Value = '0.3867'
$ws.Cells.Item(44, 5).Value = '  -1.51%  '
$ws.Cells.Item(45, 4).Value = '6.973'
$ws.Cells.Item(45, 5).Value = '  +1.61%  '
$ws.Cells.Item(46, 4).Value = '0.1162'
$ws.Cells.Item(46, 5).Value = '  -1.98%  '
$ws.Cells.Item(47, 4).Value = '0.05391'
$ws.Cells.Item(48, 4).Value = '30.53'
$ws.Cells.Item(48, 5).Value = '  -0.20%  '
$ws.Cells.Item(49, 4).Value = '7.747'
$ws.Cells.Item(49, 5).Value = '  -0.80%  '
$ws.Cells.Item(50, 4).Value = '51.80'
$ws.Cells.Item(50, 5).Value = '  -0.09%  '
$ws.Cells.Item(51, 4).Value = '1.229'
$ws.Cells.Item(51, 5).Value = '  -1.99%  '
